$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues = -4163
$xlPasteValues = -4163

$ws.Range("D2").Formula = '="69.295.12"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial($xlPasteValues)
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Formula = '="3.417.51"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial($xlPasteValues)
$ws.Range("E3").Value = "  +0.94%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Formula = '="579.49"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial($xlPasteValues)
$ws.Range("E5").Value = "  -1.57%  "

$ws.Range("D6").Formula = '="176.25"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial($xlPasteValues)
$ws.Range("E6").Value = "  -2.78%  "

$ws.Range("D8").Formula = '="3.411.58"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial($xlPasteValues)
$ws.Range("E8").Value = "  +0.91%  "

$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("D10").Formula = '="0.197"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial($xlPasteValues)
$ws.Range("E10").Value = "  +0.71%  "

$ws.Range("E11").Value = "  -1.05%  "

$ws.Range("D12").Formula = '="48.72"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial($xlPasteValues)
$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("D13").Formula = '="0.0000279"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial($xlPasteValues)
$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("D14").Formula = '="694.02"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial($xlPasteValues)
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("D15").Formula = '="3.963.63"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial($xlPasteValues)
$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("D16").Formula = '="8.62"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial($xlPasteValues)
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").Formula = '="69.348.87"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial($xlPasteValues)
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").Formula = '="3.415.97"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial($xlPasteValues)
$ws.Range("E18").Value = "  +1.20%  "

$ws.Range("E19").Value = "  +0.74%  "

$ws.Range("D20").Formula = '="17.63"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial($xlPasteValues)
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("D21").Formula = '="11.36"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial($xlPasteValues)
$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("D22").Formula = '="0.895"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial($xlPasteValues)
$ws.Range("E22").Value = "  -0.56%  "

$ws.Range("D23").Formula = '="5.42"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial($xlPasteValues)
$ws.Range("E23").Value = "  +0.34%  "

$ws.Range("D24").Formula = '="16.88"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial($xlPasteValues)
$ws.Range("E24").Value = "  -1.70%  "

$ws.Range("D25").Formula = '="100.38"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial($xlPasteValues)
$ws.Range("E25").Value = "  -4.22%  "

$ws.Range("D26").Formula = '="3.88"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial($xlPasteValues)
$ws.Range("E26").Value = "  -1.76%  "

$ws.Range("D27").Formula = '="2.66"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial($xlPasteValues)
$ws.Range("E27").Value = "  -2.61%  "

$ws.Range("D28").Formula = '="9.57"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial($xlPasteValues)
$ws.Range("E28").Value = "  -0.62%  "

$ws.Range("D29").Formula = '="33.28"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial($xlPasteValues)
$ws.Range("E29").Value = "  -3.20%  "

$ws.Range("D30").Formula = '="8.72"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial($xlPasteValues)

$ws.Range("D31").Formula = '="6.92"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial($xlPasteValues)
$ws.Range("E31").Value = "  -1.65%  "

$ws.Range("D32").Formula = '="568.86"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial($xlPasteValues)
$ws.Range("E32").Value = "  +2.01%  "

$ws.Range("D33").Formula = '="3.69"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial($xlPasteValues)
$ws.Range("E33").Value = "  +0.60%  "

$ws.Range("D34").Formula = '="10.98"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial($xlPasteValues)
$ws.Range("E34").Value = "  -1.96%  "

$ws.Range("E35").Value = "  -2.50%  "

$ws.Range("D36").Formula = '="58.16"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial($xlPasteValues)
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("D37").Formula = '="1.00"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial($xlPasteValues)
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").Formula = '="3.580.63"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial($xlPasteValues)
$ws.Range("E38").Value = "  -3.82%  "

$ws.Range("D39").Formula = '="0.138"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial($xlPasteValues)
$ws.Range("E39").Value = "  -2.49%  "

$ws.Range("D40").Formula = '="34.76"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial($xlPasteValues)
$ws.Range("E40").Value = "  -1.06%  "

$ws.Range("D41").Formula = '="0.0₃0724"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial($xlPasteValues)
$ws.Range("E41").Value = "  +2.52%  "

$ws.Range("D42").Formula = '="3.26"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial($xlPasteValues)
$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("D43").Formula = '="2.65"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial($xlPasteValues)
$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("E44").Value = "  -2.81%  "

$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("E46").Value = "  +2.64%  "

$ws.Range("E47").Value = "  -1.43%  "

$ws.Range("E48").Value = "  -1.60%  "

$ws.Range("D49").Formula = '="0.999"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial($xlPasteValues)
$ws.Range("E49").Value = "  -0.21%  "

$ws.Range("D50").Formula = '="131.83"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial($xlPasteValues)
$ws.Range("E50").Value = "  -0.57%  "

$ws.Range("D51").Formula = '="2.64"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial($xlPasteValues)
$ws.Range("E51").Value = "  +0.97%  "
